$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H94").Value = 1345.5
$ws.Range("I94").Value = 1460.6666
$ws.Range("K94").Value = 1460.6666
$ws.Range("M94").Value = -1009.6666
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()
$ws.Range("H135").Value = 5035.9165
$ws.Range("I135").Value = 4307.5
$ws.Range("J135").Value = 6492.75
$ws.Range("K135").Value = 38767.5
$ws.Range("L135").Value = 58434.75
$ws.Range("M135").Value = -36232.5
$ws.Range("N135").Value = -63504.75
$ws.Range("H137").Value = 4785.7144
$ws.Range("I137").Value = 2900.2
$ws.Range("J137").Value = 9499.5
$ws.Range("K137").Value = 8700.599999999999
$ws.Range("L137").Value = 28498.5
$ws.Range("M137").Value = -6150.599999999999
$ws.Range("N137").Value = -33598.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 1001500
$ws.Range("I10").Value = 2000000
$ws.Range("J10").Value = 3000
$ws.Range("K10").Value = 2000000
$ws.Range("L10").Value = 3000
$ws.Range("M10").Value = -1999830
$ws.Range("N10").Value = -3340
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").ClearContents()
$ws.Range("H45").Value = 9874.75
$ws.Range("I45").Value = 9874.75
$ws.Range("K45").Value = 9874.75
$ws.Range("M45").Value = -9497.75
$ws.Range("H97").Value = 1114.5454
$ws.Range("J97").Value = 1448
$ws.Range("L97").Value = 1448
$ws.Range("N97").Value = -2440
$ws.Range("H113").Value = 60928.668
$ws.Range("J113").Value = 60928.668
$ws.Range("L113").Value = 60928.668
$ws.Range("N113").Value = -69606.66800000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 175
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 175
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 175
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = -401
$ws.Range("H86").Value = 1392.25
$ws.Range("I86").Value = 1234.75
$ws.Range("J86").Value = 1549.75
$ws.Range("K86").Value = 1234.75
$ws.Range("L86").Value = 1549.75
$ws.Range("M86").Value = -111.75
$ws.Range("N86").Value = -3795.75
$ws.Range("H89").Value = 1392.25
$ws.Range("I89").Value = 1234.75
$ws.Range("J89").Value = 1549.75
$ws.Range("K89").Value = 6173.75
$ws.Range("L89").Value = 7748.75
$ws.Range("M89").Value = -557.75
$ws.Range("N89").Value = -18980.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 23.75
$ws.Range("I19").Value = 23.75
$ws.Range("K19").Value = 23.75
$ws.Range("M19").Value = 146.25
$ws.Range("H24").Value = 23.75
$ws.Range("I24").Value = 23.75
$ws.Range("K24").Value = 23.75
$ws.Range("M24").Value = 146.25
$ws.Range("H31").Value = 2248.0588
$ws.Range("I31").Value = 1915.909
$ws.Range("J31").Value = 2857
$ws.Range("K31").Value = 1915.909
$ws.Range("L31").Value = 2857
$ws.Range("M31").Value = -1620.909
$ws.Range("N31").Value = -3447
$ws.Range("H34").Value = 2248.0588
$ws.Range("I34").Value = 1915.909
$ws.Range("J34").Value = 2857
$ws.Range("K34").Value = 1915.909
$ws.Range("L34").Value = 2857
$ws.Range("M34").Value = -1713.909
$ws.Range("N34").Value = -3261
$ws.Range("H99").Value = 1349.1538
$ws.Range("I99").Value = 999
$ws.Range("K99").Value = 999
$ws.Range("M99").Value = 499
$ws.Range("H100").Value = 175390
$ws.Range("J100").Value = 175390
$ws.Range("L100").Value = 175390
$ws.Range("N100").Value = -177554
$ws.Range("H105").Value = 3418.8462
$ws.Range("I105").Value = 2489
$ws.Range("K105").Value = 2489
$ws.Range("M105").Value = -742
$ws.Range("H122").Value = 1490.8334
$ws.Range("J122").Value = 2200
$ws.Range("L122").Value = 6600
$ws.Range("N122").Value = -11500
$ws.Range("H126").Value = 1349.1538
$ws.Range("I126").Value = 999
$ws.Range("K126").Value = 2997
$ws.Range("M126").Value = -527
$ws.Range("H141").Value = 40132.4
$ws.Range("J141").Value = 40132.4
$ws.Range("L141").Value = 40132.4
$ws.Range("N141").Value = -50492.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("M43").ClearContents()
$ws.Range("N43").ClearContents()
$ws.Range("H44").Value = 488.9091
$ws.Range("I44").Value = 274
$ws.Range("K44").Value = 822
$ws.Range("M44").Value = -424
$ws.Range("H55").Value = 5661.6665
$ws.Range("J55").Value = 5661.6665
$ws.Range("L55").Value = 16984.9995
$ws.Range("N55").Value = -17338.9995
$ws.Range("H60").Value = 566
$ws.Range("I60").Value = 566
$ws.Range("K60").Value = 1698
$ws.Range("M60").Value = -1447
$ws.Range("H97").Value = 518.75
$ws.Range("J97").Value = 548
$ws.Range("L97").Value = 1644
$ws.Range("N97").Value = -2636

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 30000
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 30000
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 30000
$ws.Range("M18").ClearContents()
$ws.Range("N18").Value = -30586
$ws.Range("H20").Value = 450000
$ws.Range("I20").Value = 450000
$ws.Range("K20").Value = 450000
$ws.Range("M20").Value = -449755
$ws.Range("H97").Value = 283.33334
$ws.Range("I97").Value = 100
$ws.Range("K97").Value = 100
$ws.Range("M97").Value = 396
$ws.Range("H113").Value = 999
$ws.Range("I113").Value = 999
$ws.Range("K113").Value = 999
$ws.Range("M113").Value = 1171
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()
$ws.Range("H122").Value = 3035.9412
$ws.Range("I122").Value = 3107.4
$ws.Range("K122").Value = 9322.200000000001
$ws.Range("M122").Value = -6872.200000000001
$ws.Range("H134").Value = 36475
$ws.Range("J134").Value = 36475
$ws.Range("L134").Value = 109425
$ws.Range("N134").Value = -114495

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 9512.733
$ws.Range("I7").Value = 10225
$ws.Range("J7").Value = 9253.727999999999
$ws.Range("K7").Value = 10225
$ws.Range("L7").Value = 9253.727999999999
$ws.Range("M7").Value = -10113
$ws.Range("N7").Value = -9477.727999999999
$ws.Range("H43").Value = 34900
$ws.Range("I43").Value = 34900
$ws.Range("K43").Value = 34900
$ws.Range("M43").Value = -34707
$ws.Range("H61").Value = 1899.6666
$ws.Range("I61").Value = 1899.6666
$ws.Range("K61").Value = 1899.6666
$ws.Range("M61").Value = -1697.6666
$ws.Range("H68").Value = 4000
$ws.Range("I68").Value = 4000
$ws.Range("K68").Value = 4000
$ws.Range("M68").Value = -3251
$ws.Range("H71").Value = 4000
$ws.Range("I71").Value = 4000
$ws.Range("K71").Value = 20000
$ws.Range("M71").Value = -16256
$ws.Range("H113").Value = 1899.6666
$ws.Range("I113").Value = 1899.6666
$ws.Range("K113").Value = 1899.6666
$ws.Range("M113").Value = 270.3334
$ws.Range("H122").Value = 7801.269
$ws.Range("I122").Value = 7970.8335
$ws.Range("J122").Value = 7655.9287
$ws.Range("K122").Value = 23912.5005
$ws.Range("L122").Value = 22967.7861
$ws.Range("M122").Value = -21462.5005
$ws.Range("N122").Value = -27867.7861
$ws.Range("H126").Value = 9512.733
$ws.Range("I126").Value = 10225
$ws.Range("J126").Value = 9253.727999999999
$ws.Range("K126").Value = 30675
$ws.Range("L126").Value = 27761.184
$ws.Range("M126").Value = -28205
$ws.Range("N126").Value = -32701.184
$ws.Range("H136").Value = 1214.2858
$ws.Range("I136").Value = 1214.2858
$ws.Range("K136").Value = 3642.8574
$ws.Range("M136").Value = -1092.8574

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1763.6364
$ws.Range("I96").Value = 2001
$ws.Range("J96").Value = 1628
$ws.Range("K96").Value = 2001
$ws.Range("L96").Value = 1628
$ws.Range("M96").Value = -628
$ws.Range("N96").Value = -4374
$ws.Range("H113").Value = 614.3333
$ws.Range("J113").Value = 671.5
$ws.Range("L113").Value = 2014.5
$ws.Range("N113").Value = -6354.5
$ws.Range("H122").Value = 1182.1428
$ws.Range("I122").Value = 712.5
$ws.Range("K122").Value = 2137.5
$ws.Range("M122").Value = 312.5
$ws.Range("H126").Value = 2985.1428
$ws.Range("I126").Value = 2499.2
$ws.Range("K126").Value = 7497.599999999999
$ws.Range("M126").Value = -5027.599999999999
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()
$ws.Range("H132").Value = 7000
$ws.Range("I132").Value = 7000
$ws.Range("K132").Value = 21000
$ws.Range("M132").Value = -18470
